$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-02"

# Update the column header label (shared string) for the 2022 total column
$ws.Range("I1").Value = "2022 (through 06-02)"

# Update July row (row 7): 2021 total (H7) and 2022-through-date total (I7)
$ws.Range("H7").Value = 129
$ws.Range("I7").Value = 6

# Update the grand Total row (row 14): 2021 total (H14) and 2022-through-date total (I14)
$ws.Range("H14").Value = 1849
$ws.Range("I14").Value = 670
